# Fix list level numbering: the content of a top-level list should be
# at the same indent level as a top-level paragraph (only continuation
# paragraphs of a list should be nested one level deeper).
#
# Note: TextRange.Paragraphs(...).IndentLevel is 1-based (OOXML lvl="0"
# corresponds to IndentLevel 1), so decrementing IndentLevel by one here
# turns the OOXML lvl="N" into lvl="N-1".

$p = $ppt.ActivePresentation

# Slide 1, Content Placeholder 2:
#   "Bullet item with inline code" : lvl 1 -> 0
#   "with nested"                 : lvl 2 -> 1
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$tr1.Paragraphs(2, 1).IndentLevel = 1
$tr1.Paragraphs(4, 1).IndentLevel = 2

# Slide 2, Content Placeholder 2:
#   "Nested" : lvl 1 -> 0
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2, 1).IndentLevel = 1

# Slide 3, Content Placeholder 2:
#   "A total alternative for head" : lvl 1 -> 0
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Paragraphs(1, 1).IndentLevel = 1
